$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Move the active selection to C38 (matches the saved view's sqref/activeCell)
$ws.Range("C38").Select()

# Widen the "Designations" column (D) so its contents fit without wrapping
$ws.Columns.Item(4).ColumnWidth = 40.166666666666664

# Hide the "Manufacturer" (G) and "Manufacturer Part #" (H) columns
$ws.Columns.Item(7).Hidden = $true
$ws.Columns.Item(8).Hidden = $true

# Shrink the print scale from 59% to 72%
$ws.PageSetup.Zoom = 72
